# Adds two new configuration rows ("PECO" and "BGE") that duplicate the
# existing "default" (imperial) config row (row 2), appended below the
# current data as rows 4 and 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column values (B..AC) shared by the existing "default" imperial config
# row (row 2) that the new PECO/BGE rows duplicate verbatim; only column A
# (the config name) differs between the new rows.
$commonValues = @(
    0.9, 0.7, 0, "ft/s", 2.533, "ft/s", 90, "degrees", 90,
    "degrees", 200, "ft", "N-S", 40, 10, 6, 2009, 1400, "imperial",
    "Industrial", -15, 40, 5, 35, "C", 35, 10, 15
)

$names = @("PECO", "BGE")

for ($i = 0; $i -lt $names.Length; $i++) {
    $targetRow = 4 + $i

    $arr = New-Object 'object[,]' 1, (1 + $commonValues.Length)
    $arr[0, 0] = $names[$i]
    for ($c = 0; $c -lt $commonValues.Length; $c++) {
        $arr[0, $c + 1] = $commonValues[$c]
    }

    $startCell = $ws.Cells.Item($targetRow, 1)
    $endCell = $ws.Cells.Item($targetRow, 1 + $commonValues.Length)
    $rng = $ws.Range($startCell, $endCell)
    $rng.Value = $arr
}

# Move the active selection to A6, matching the post-edit state.
$ws.Range("A6").Select()
